$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.18%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.73%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.039"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.07%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07466"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.79%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.355"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.02%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.75%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9276"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'2.424"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.13%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.27%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1826"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.11%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08863"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.75%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04165"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.12%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001278"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.55%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005933"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.44%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.344"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.92%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3291"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.23%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'7.881"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.68%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1409"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.76%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2966"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.80%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04039"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.75%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001264"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.57%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.003878"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.02%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-4.26%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.39%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02393"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'4.22%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'5.75%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006879"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-3.50%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.33%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1321"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.34%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007375"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.33%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007168"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.44%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3216"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.49%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006232"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.98%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.38%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.04601"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-81.73%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004200"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.38%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.38%"
$ws.Range("E51").Style = "Normal"
